$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 100, shifting existing rows 100-193 down to 101-194
$ws.Rows.Item(100).Insert()

# Populate the new row 100 with the new record
$ws.Cells.Item(100, 1).Value = 6
$ws.Cells.Item(100, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(100, 3).Value = "Metropolitana"
$ws.Cells.Item(100, 4).Value = 44658
$ws.Cells.Item(100, 5).Value = 13
$ws.Cells.Item(100, 6).Value = "Fruta"
$ws.Cells.Item(100, 7).Value = 100101
$ws.Cells.Item(100, 8).Value = "Berries"
$ws.Cells.Item(100, 9).Value = 100101004
$ws.Cells.Item(100, 10).Value = "Frambuesa"
$ws.Cells.Item(100, 11).Value = "Sin especificar"
$ws.Cells.Item(100, 12).Value = "Primera"
$ws.Cells.Item(100, 13).Value = 350
$ws.Cells.Item(100, 14).Value = 8000
$ws.Cells.Item(100, 15).Value = 8000
$ws.Cells.Item(100, 16).Value = 8000
$ws.Cells.Item(100, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(100, 18).Value = "Provincia de Linares"
$ws.Cells.Item(100, 19).Value = 4000
$ws.Cells.Item(100, 20).Value = 2
